# Switch to platform IDs for stop IDs
# Row 3 (D3:Y3) currently holds the GTFS stop "name" slugs (e.g. "san_francisco")
# as text; replace them with the corresponding numeric platform stop IDs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$platformIds = @(70012, 70022, 70032, 70042, 70052, 70062, 70082, 70092, 70102, 70112, 70122, 70132, 70142, 70162, 70172, 70192, 70202, 70212, 70222, 70232, 70242, 70262)

$col = 4  # column D
foreach ($id in $platformIds) {
    $ws.Cells.Item(3, $col).Value = $id
    $col = $col + 1
}
